$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

# --- Simple price (column D) updates ---
Set-TextValue "D2" "245.35"
Set-TextValue "D3" "22.05"
Set-TextValue "D5" "0.05853"
Set-TextValue "D6" "3.394"
Set-TextValue "D8" "0.8133"
Set-TextValue "D9" "1.017"

# --- Rows 10-18: coin listing rotates up by one (wraps around) ---
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D10" "0.1422"
$ws.Range("E10").Value = "9WazirXWRX"

$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D11" "0.04186"
$ws.Range("E11").Value = "10LiechtensteinCryptoassetsExchangeLCXBestin24h"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D12" "0.07394"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D13" "0.02984"
$ws.Range("E13").Value = "12BitrueCoinBTR"

$ws.Range("B14").Value = "MCDex"
$ws.Range("C14").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D14" "4.136"
$ws.Range("E14").Value = "13MCDexMCB"

$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D15" "0.09392"
$ws.Range("E15").Value = "14BitMartTokenBMX"

$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D16" "0.001586"
$ws.Range("E16").Value = "15BitForexTokenBF"

$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D17" "0.04828"
$ws.Range("E17").Value = "16CoinExTokenCET"

$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D18" "0.0005890"
$ws.Range("E18").Value = "17OneONE"

# --- More simple price updates ---
Set-TextValue "D19" "0.005937"
Set-TextValue "D20" "0.004081"
Set-TextValue "D21" "0.0009877"
Set-TextValue "D23" "3.715"
Set-TextValue "D24" "2.227"
Set-TextValue "D27" "0.0002483"
Set-TextValue "D40" "0.03866"

# --- Rows 41-43: coin listing rotates up by one (wraps around) ---
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D41" "0.1073"
$ws.Range("E41").Value = "40BKEXTokenBKK"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.002412"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D43" "0.003023"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"

# --- Row 44: price update + label text change (rank string only, no rotation) ---
Set-TextValue "D44" "0.005070"
$ws.Range("E44").Value = "43LocalTradersLCT"

# --- Final simple price updates ---
Set-TextValue "D47" "0.7700"
Set-TextValue "D49" "0.00002100"
